$wb = $excel.ActiveWorkbook

# ---- Sheet "Pasos": fill in the new step/expected-result text (rows 3-5) ----
$wsPasos = $wb.Worksheets.Item("Pasos")
$wsPasos.Range("B3").Value = "Selecciono la opcion Ingresar"
$wsPasos.Range("C3").Value = "El sistema muestra un popup para iniciar sesion o registrarse."

# ---- Sheet "Precondiciones": precondition text updated (usuario ya existe) ----
$wsPrecond = $wb.Worksheets.Item("Precondiciones")
$wsPrecond.Range("B2").Value = "El usuario ""CPA_Usuario"" existe en la base de datos"
$wsPrecond.Rows.Item(2).EntireRow.AutoFit()

# ---- back to "Pasos": remaining new rows ----
$wsPasos.Range("B4").Value = "Ingreso ""CPA_Usuario"" en el campo usuario y ""CPA_Contraseña"" en el campo contraseña"
$wsPasos.Range("B5").Value = "Seleccióno la opcion Iniciar sesion"
$wsPasos.Range("C5").Value = "El sistema muestra el nombre de usuario donde antes estaba la opcion Ingresar"

# Rows 3-5 grew to two lines of wrapped text
$wsPasos.Rows.Item(3).RowHeight = 26.25
$wsPasos.Rows.Item(4).RowHeight = 26.25
$wsPasos.Rows.Item(5).RowHeight = 26.25

# ---- selections (cosmetic cursor position per sheet) ----
$wsGenerales = $wb.Worksheets.Item("DatosGenerales")
$wsGenerales.Range("B24").Select()

$wsPrecond.Range("B4").Select()

# "Pasos" stays the active/tabSelected sheet, so select it last
$wsPasos.Range("B6").Select()
